# Update leve profit data across multiple sheets (scheduled data refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 310
$ws.Range("I2").Value = 265
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 265
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -152
$ws.Range("N2").Value = -626
$ws.Range("H138").Value = 2901196.5
$ws.Range("I138").Value = 1089.6552
$ws.Range("J138").Value = 5003774
$ws.Range("K138").Value = 3268.9656
$ws.Range("L138").Value = 15011322
$ws.Range("M138").Value = 1871.0344
$ws.Range("N138").Value = -15021602
$ws.Range("H141").Value = 1651.8889
$ws.Range("I141").Value = 884.9091
$ws.Range("J141").Value = 2857.1428
$ws.Range("K141").Value = 2654.7273
$ws.Range("L141").Value = 8571.428400000001
$ws.Range("M141").Value = 2525.2727
$ws.Range("N141").Value = -18931.4284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 712.9375
$ws.Range("I97").Value = 515.9231
$ws.Range("J97").Value = 1566.6666
$ws.Range("K97").Value = 515.9231
$ws.Range("L97").Value = 1566.6666
$ws.Range("M97").Value = -19.92309999999998
$ws.Range("N97").Value = -2558.6666
$ws.Range("H132").Value = 3310.4285
$ws.Range("I132").Value = 3061.2964
$ws.Range("J132").Value = 4151.25
$ws.Range("K132").Value = 9183.889200000001
$ws.Range("L132").Value = 12453.75
$ws.Range("M132").Value = -6653.889200000001
$ws.Range("N132").Value = -17513.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 31200
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 31200
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 31200
$ws.Range("N117").Value = -40378
$ws.Range("H118").Value = 33542.855
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 33542.855
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 33542.855
$ws.Range("N118").Value = -36856.855
$ws.Range("H119").Value = 30000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 30000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws.Range("H120").Value = 30000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 30000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H122").Value = 20000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 20000
$ws.Range("N122").Value = -29800
$ws.Range("H123").Value = 25000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 25000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800
$ws.Range("H124").Value = 32190
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 32190
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 32190
$ws.Range("N124").Value = -42010
$ws.Range("H125").Value = 20000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 20000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -29840
$ws.Range("H126").Value = 33780
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 33780
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 33780
$ws.Range("N126").Value = -43660
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 3000
$ws.Range("I128").Value = 3000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 9000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -6510
$ws.Range("H129").Value = 35559.082
$ws.Range("I129").Value = 30709
$ws.Range("J129").Value = 36000
$ws.Range("K129").Value = 30709
$ws.Range("L129").Value = 36000
$ws.Range("M129").Value = -25709
$ws.Range("N129").Value = -46000
$ws.Range("H130").Value = 118000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 118000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 118000
$ws.Range("N130").Value = -128040
$ws.Range("H131").Value = 30000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H133").Value = 42000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 42000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 42000
$ws.Range("N133").Value = -52120
$ws.Range("H134").Value = 3787.5881
$ws.Range("I134").Value = 2588.762
$ws.Range("J134").Value = 5724.154
$ws.Range("K134").Value = 7766.286
$ws.Range("L134").Value = 17172.462
$ws.Range("M134").Value = -5231.286
$ws.Range("N134").Value = -22242.462
$ws.Range("H135").Value = 45750
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 45750
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 45750
$ws.Range("N135").Value = -55890
$ws.Range("H137").Value = 39939.395
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 39939.395
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 39939.395
$ws.Range("N137").Value = -50139.395
$ws.Range("H138").Value = 38770
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 38770
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 38770
$ws.Range("N138").Value = -49050
$ws.Range("H139").Value = 49214.285
$ws.Range("I139").Value = 54000
$ws.Range("J139").Value = 48846.152
$ws.Range("K139").Value = 54000
$ws.Range("L139").Value = 48846.152
$ws.Range("M139").Value = -48860
$ws.Range("N139").Value = -59126.152
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 99000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 99000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 99000
$ws.Range("N141").Value = -109360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8306.177
$ws.Range("J68").Value = 2854.818
$ws.Range("L68").Value = 2854.818
$ws.Range("N68").Value = -4352.818
$ws.Range("H71").Value = 8306.177
$ws.Range("J71").Value = 2854.818
$ws.Range("L71").Value = 14274.09
$ws.Range("N71").Value = -21762.09

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3854.0312
$ws.Range("J107").Value = 622.5
$ws.Range("L107").Value = 1867.5
$ws.Range("N107").Value = -5707.5
$ws.Range("H132").Value = 2901.617
$ws.Range("I132").Value = 2829.45
$ws.Range("J132").Value = 3314
$ws.Range("K132").Value = 8488.349999999999
$ws.Range("L132").Value = 9942
$ws.Range("M132").Value = -5958.349999999999
$ws.Range("N132").Value = -15002

Write-Host "All changes applied"
